# Añadir formulario de informe en la sección 1
# Extend the "tabla_pf" report table (Sheet1) from 2 data rows (iterations 0-1)
# to 10 data rows (iterations 0-9), updating the numeric columns xn/fxn/E with
# the full iteration series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Iteration, xn, fxn, E  (row 1 header "Iteration/xn/fxn/E" is unchanged)
$data = @(
    @("0", "0.1",                "0.0701665833531719", "1.0005"),
    @("1", "0.0649167083234141", "0.0523197567634787", "0.0350832916765859"),
    @("2", "0.0387568299416747", "0.0342602563314401", "0.0261598783817393"),
    @("3", "0.0216267017759547", "0.0202252449001057", "0.0171301281657201"),
    @("4", "0.0115140793259018", "0.0111166116673483", "0.0101126224500529"),
    @("5", "0.0059557734922276", "0.0058493949882694", "0.0055583058336741"),
    @("6", "0.0030310759980929", "0.0030035183742677", "0.0029246974941347"),
    @("7", "0.001529316810959",  "0.0015223009773643", "0.0015017591871338"),
    @("8", "0.0007681663222768", "0.0007663961593273", "0.0007611504886821"),
    @("9", "0.0003849682426131", "0.0003845236504784", "0.0003831980796636")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    # Force text entry (leading apostrophe) so the decimal strings are stored
    # verbatim instead of being coerced/rounded as floating point numbers --
    # matching the source table's existing text-typed cells.
    $ws.Cells.Item($row, 1).Value = "'" + $data[$i][0]
    $ws.Cells.Item($row, 2).Value = "'" + $data[$i][1]
    $ws.Cells.Item($row, 3).Value = "'" + $data[$i][2]
    $ws.Cells.Item($row, 4).Value = "'" + $data[$i][3]

    # The apostrophe-forced text entry marks the cell with the "number
    # stored as text" quote-prefix style; reset back to the plain/default
    # style so the data rows keep the workbook's original (unstyled) look.
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 3).Style = "Normal"
    $ws.Cells.Item($row, 4).Style = "Normal"
}
